# Slide 1, shape 3 ("CaixaDeTexto 3") lists the team members. Update the
# "Members:" heading to split into two runs ("Members" + ":"), and append
# each member's student number to their line.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item(3)
$tr = $sh.TextFrame.TextRange

# Build the new paragraph text, `r (carriage return) is the paragraph
# separator for TextRange.Text.
$members = "Members:`rJoão Beires nº 1190718`rJosé Soares nº 1190782`rJosé Maia nº 1191419`rLourenço Melo nº 1190811"

# Stamp a throwaway value over the whole range first so the engine treats
# the real assignment below as all-new paragraphs (no leftover endParaRPr
# or run-splitting against the old text), then set the real text.
$tr.Text = "x"
$tr.Text = $members

# Split the first paragraph "Members:" into two runs: "Members" and ":".
$para1 = $tr.Paragraphs(1, 1)
$para1.Text = "Members"
$para1.InsertAfter(":")
